$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the species/record data between row 2 and row 3 ---
# Column A (Id)
$a2 = $ws.Range("A2").Value2
$a3 = $ws.Range("A3").Value2
$ws.Range("A2").Value = $a3
$ws.Range("A3").Value = $a2

# Column B (Taxonsorteringsordning)
$b2 = $ws.Range("B2").Value2
$b3 = $ws.Range("B3").Value2
$ws.Range("B2").Value = $b3
$ws.Range("B3").Value = $b2

# Column D (Rödlistade)
$d2 = $ws.Range("D2").Value2
$d3 = $ws.Range("D3").Value2
$ws.Range("D2").Value = $d3
$ws.Range("D3").Value = $d2

# Column E (TaxonId)
$e2 = $ws.Range("E2").Value2
$e3 = $ws.Range("E3").Value2
$ws.Range("E2").Value = $e3
$ws.Range("E3").Value = $e2

# Column F (Artnamn)
$f2 = $ws.Range("F2").Value2
$f3 = $ws.Range("F3").Value2
$ws.Range("F2").Value = $f3
$ws.Range("F3").Value = $f2

# Column G (Vetenskapligt namn)
$g2 = $ws.Range("G2").Value2
$g3 = $ws.Range("G3").Value2
$ws.Range("G2").Value = $g3
$ws.Range("G3").Value = $g2

# Column H (Auktor)
$h2 = $ws.Range("H2").Value2
$h3 = $ws.Range("H3").Value2
$ws.Range("H2").Value = $h3
$ws.Range("H3").Value = $h2

# --- Update coordinates (Ost / Nord) with rounded, recalculated values ---
$ws.Range("Q2").Value = 792456
$ws.Range("R2").Value = 7344877
$ws.Range("Q3").Value = 792464
$ws.Range("R3").Value = 7344871

# --- Clear the Starttid / Sluttid columns for both rows ---
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
$ws.Range("Z3").ClearContents()
$ws.Range("AB3").ClearContents()
